$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card1")

# Fix header N1: remove trailing space from "Correction "
$ws.Range("N1").Value = "Correction"

# Add new header O1: "Serviced by " (trailing space intentional)
# Copy N1's formatting (bold, borders, alignment) onto O1 first, then set its text
$ws.Range("N1").Copy($ws.Range("O1"))
$ws.Range("O1").Value = "Serviced by "

# Fill N2:N12 with "nan" placeholder text (matches other "nan" columns)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}

# Create O2:O12 as present-but-empty text cells (mirrors N2:N12's pre-edit
# blank-placeholder shape). A bare Value="" is treated as "no entry" and the
# cell disappears again, so write the text-prefix marker "'" (quote-only
# input collapses to an empty string) then strip the formatting it implies.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 15).Value = "'"
    $ws.Cells.Item($r, 15).ClearFormats()
}
